$wb = $excel.ActiveWorkbook

# --- "actual_cost_v1" sheet: fill in missing Quantity (col B) and GrandTotal (col E) cells ---
$ws4 = $wb.Worksheets.Item("actual_cost_v1")
$ws4.Range("B2").Value = 0
$ws4.Range("E3").Value = 2869
$ws4.Range("E4").Value = 2869
$ws4.Range("B5").Value = 0
$ws4.Range("E5").Value = 2869
$ws4.Range("B6").Value = 0
$ws4.Range("E6").Value = 2869
$ws4.Range("B7").Value = 0
$ws4.Range("E7").Value = 2869
$ws4.Range("B8").Value = 0
$ws4.Range("E8").Value = 2869
$ws4.Range("B9").Value = 0
$ws4.Range("E9").Value = 2869
$ws4.Range("B10").Value = 0
$ws4.Range("E10").Value = 2869
$ws4.Range("B11").Value = 0
$ws4.Range("E11").Value = 2869
$ws4.Range("B12").Value = 0
$ws4.Range("E12").Value = 2869
$ws4.Range("B13").Value = 0
$ws4.Range("E13").Value = 2869
$ws4.Range("E14").Value = 2869
$ws4.Range("B15").Value = 0
$ws4.Range("E15").Value = 2869
$ws4.Range("E16").Value = 2869
$ws4.Range("E17").Value = 2869
$ws4.Range("B18").Value = 0
$ws4.Range("E18").Value = 2869

# Update this sheet's remembered selection (does not switch the active tab
# as long as a later Select() on another sheet happens afterwards).
$ws4.Range("L16").Select()

# --- "planned_estimated_cost_v1" sheet: fill in missing GrandTotal (col E) cells ---
$ws6 = $wb.Worksheets.Item("planned_estimated_cost_v1")
$ws6.Range("E3").Value = 23784
$ws6.Range("E4").Value = 23784
$ws6.Range("E5").Value = 23784
$ws6.Range("E6").Value = 23784
$ws6.Range("E7").Value = 23784
$ws6.Range("E8").Value = 23784
$ws6.Range("E9").Value = 23784
$ws6.Range("E10").Value = 23784
$ws6.Range("E11").Value = 23784
$ws6.Range("E12").Value = 23784
$ws6.Range("E13").Value = 23784
$ws6.Range("E14").Value = 23784
$ws6.Range("E15").Value = 23784
$ws6.Range("E16").Value = 23784
$ws6.Range("E17").Value = 23784

# Make this the active sheet/tab, with the newly filled column selected
# (this is the sheet the user ended the session on).
$ws6.Range("E2:E17").Select()
